# Enable parallel test execution and add logging:
# append new order IDs to the "Orders" sheet, matching the newly
# generated / logged test orders.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

$newOrderIds = @(
    "000054816",
    "000054817",
    "000054818",
    "000054819",
    "000054865",
    "000054866",
    "000054890",
    "000054891"
)

$startRow = 5
for ($i = 0; $i -lt $newOrderIds.Count; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    # Enter as a text formula so Excel keeps the leading zeros, then
    # freeze it into a plain (shared-string) value via copy / paste-special
    # so no number-format / style gets attached to the cell.
    $cell.Formula = "=""" + $newOrderIds[$i] + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
